$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("MainSheet")
$ws2 = $wb.Worksheets.Item("ClearanceRules")

try {
    $ws1.Activate()
    Write-Output "activate ok"
} catch {
    Write-Output "activate failed: $_"
}

try {
    $ws1.Range("K7").Select() | Out-Null
    Write-Output "select ok"
} catch {
    Write-Output "select failed: $_"
}
